$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '22.392.84'
Set-TextValue 2 5 '  -0.13%  '
Set-TextValue 3 4 '1.566.02'
Set-TextValue 3 5 '  -0.06%  '
Set-TextValue 4 5 '  -0.09%  '
Set-TextValue 5 5 '  -0.13%  '
Set-TextValue 6 4 '286.11'
Set-TextValue 6 5 '  +0.34%  '
Set-TextValue 7 4 '0.3713'
Set-TextValue 7 5 '  +2.40%  '
Set-TextValue 8 4 '0.3275'
Set-TextValue 9 4 '46.49'
Set-TextValue 9 5 '  -4.09%  '
Set-TextValue 10 4 '1.144'
Set-TextValue 10 5 '  +2.02%  '
Set-TextValue 11 4 '0.07410'
Set-TextValue 11 5 '  +0.27%  '
Set-TextValue 12 5 '  -0.05%  '
Set-TextValue 13 4 '20.42'
Set-TextValue 13 5 '  -1.61%  '
Set-TextValue 14 4 '5.837'
Set-TextValue 14 5 '  -1.73%  '
Set-TextValue 15 4 '6.815'
Set-TextValue 15 5 '  -1.18%  '
Set-TextValue 16 4 '1.566.51'
Set-TextValue 16 5 '  -0.10%  '
Set-TextValue 17 4 '0.00001097'
Set-TextValue 17 5 '  -0.53%  '
Set-TextValue 18 4 '0.06694'
Set-TextValue 18 5 '  -0.06%  '
Set-TextValue 19 4 '86.01'
Set-TextValue 19 5 '  -2.20%  '
Set-TextValue 20 4 '0.9997'
Set-TextValue 20 5 '  -0.22%  '
Set-TextValue 21 4 '6.321'
Set-TextValue 21 5 '  -0.07%  '
Set-TextValue 22 4 '16.23'
Set-TextValue 22 5 '  +0.28%  '
Set-TextValue 23 4 '11.76'
Set-TextValue 23 5 '  -1.99%  '
Set-TextValue 24 4 '22.391.99'
Set-TextValue 24 5 '  -0.11%  '
Set-TextValue 25 4 '2.300'
Set-TextValue 25 5 '  -3.14%  '
Set-TextValue 26 4 '2.559'
Set-TextValue 26 5 '  +0.88%  '
Set-TextValue 27 4 '150.94'
Set-TextValue 27 5 '  +0.36%  '
Set-TextValue 28 4 '19.32'
Set-TextValue 29 4 '4.942'
Set-TextValue 29 5 '  -1.02%  '
Set-TextValue 30 4 '123.60'
Set-TextValue 30 5 '  -0.10%  '
Set-TextValue 31 4 '1.742.94'
Set-TextValue 31 5 '  -0.05%  '
Set-TextValue 32 4 '1.046'
Set-TextValue 32 5 '  +0.96%  '
Set-TextValue 33 4 '1.956'
Set-TextValue 33 5 '  -2.43%  '
Set-TextValue 34 4 '5.934'
Set-TextValue 34 5 '  -2.49%  '
Set-TextValue 35 4 '9.629'
Set-TextValue 35 5 '  -1.77%  '
Set-TextValue 36 5 '  -0.37%  '
Set-TextValue 37 4 '1.316'
Set-TextValue 37 5 '  +2.31%  '
Set-TextValue 38 4 '0.02375'
Set-TextValue 38 5 '  -1.41%  '
Set-TextValue 39 4 '0.06286'
Set-TextValue 39 5 '  -2.12%  '
Set-TextValue 40 4 '0.2183'
Set-TextValue 40 5 '  -2.22%  '
Set-TextValue 41 4 '5.230'
Set-TextValue 41 5 '  -2.46%  '
Set-TextValue 42 4 '11.10'
Set-TextValue 42 5 '  -0.62%  '
Set-TextValue 43 4 '0.6103'
Set-TextValue 43 5 '  -2.42%  '
Set-TextValue 44 4 '1.001'
Set-TextValue 44 5 '  -0.11%  '
Set-TextValue 45 4 '13.71'
Set-TextValue 45 5 '  -0.20%  '
Set-TextValue 46 4 '0.5940'
Set-TextValue 46 5 '  -1.67%  '
Set-TextValue 47 4 '3.743'
Set-TextValue 47 5 '  -0.20%  '
Set-TextValue 48 4 '2.003'
Set-TextValue 48 5 '  -1.16%  '
Set-TextValue 49 4 '123.49'
Set-TextValue 49 5 '  +0.18%  '
Set-TextValue 50 4 '1.178'
Set-TextValue 50 5 '  -2.62%  '
Set-TextValue 51 4 '0.07151'
Set-TextValue 51 5 '  -0.66%  '
